$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Cells.Item(2, 8).Value = 2515.158
$ws.Cells.Item(2, 9).Value = 3853.9092
$ws.Cells.Item(2, 10).Value = 674.375
$ws.Cells.Item(2, 11).Value = 3853.9092
$ws.Cells.Item(2, 12).Value = 674.375
$ws.Cells.Item(2, 13).Value = -3740.9092
$ws.Cells.Item(2, 14).Value = -900.375

# Row 43
$ws.Cells.Item(43, 8).Value = 4867.222
$ws.Cells.Item(43, 9).Value = 1826
$ws.Cells.Item(43, 10).Value = 7300.2
$ws.Cells.Item(43, 11).Value = 1826
$ws.Cells.Item(43, 12).Value = 7300.2
$ws.Cells.Item(43, 14).Value = -7438.2
$ws.Cells.Item(43, 13).Value = -1757

# Row 70
$ws.Cells.Item(70, 8).Value = 7991.9644

# Row 73
$ws.Cells.Item(73, 8).Value = 7991.9644

# Row 76
$ws.Cells.Item(76, 8).Value = 5751.0713
$ws.Cells.Item(76, 9).Value = 4668
$ws.Cells.Item(76, 10).Value = 6184.3
$ws.Cells.Item(76, 11).Value = 4668
$ws.Cells.Item(76, 12).Value = 6184.3
$ws.Cells.Item(76, 13).Value = -4353
$ws.Cells.Item(76, 14).Value = -6814.3

# Row 79
$ws.Cells.Item(79, 8).Value = 5751.0713
$ws.Cells.Item(79, 9).Value = 4668
$ws.Cells.Item(79, 10).Value = 6184.3
$ws.Cells.Item(79, 11).Value = 4668
$ws.Cells.Item(79, 12).Value = 6184.3
$ws.Cells.Item(79, 13).Value = -3576
$ws.Cells.Item(79, 14).Value = -8368.299999999999

# Row 82
$ws.Cells.Item(82, 8).Value = 1199
$ws.Cells.Item(82, 10).Value = 1199
$ws.Cells.Item(82, 12).Value = 3597
$ws.Cells.Item(82, 14).Value = -4409

# Row 85
$ws.Cells.Item(85, 8).Value = 1199
$ws.Cells.Item(85, 10).Value = 1199
$ws.Cells.Item(85, 12).Value = 3597
$ws.Cells.Item(85, 14).Value = -6405

# Row 129
$ws.Cells.Item(129, 14).ClearContents()
$ws.Cells.Item(129, 8).Value = 62500936
$ws.Cells.Item(129, 9).Value = 62500936
$ws.Cells.Item(129, 10).Value = 0
$ws.Cells.Item(129, 11).Value = 187502808
$ws.Cells.Item(129, 12).Value = 0
$ws.Cells.Item(129, 13).Value = -187497808

# Row 132
$ws.Cells.Item(132, 8).Value = 23258508
$ws.Cells.Item(132, 9).Value = 23812254
$ws.Cells.Item(132, 11).Value = 71436762
$ws.Cells.Item(132, 13).Value = -71434232

# Row 138
$ws.Cells.Item(138, 8).Value = 3318.3691
$ws.Cells.Item(138, 9).Value = 2426.7273
$ws.Cells.Item(138, 10).Value = 3500
$ws.Cells.Item(138, 11).Value = 7280.1819
$ws.Cells.Item(138, 12).Value = 10500
$ws.Cells.Item(138, 13).Value = -2140.1819
$ws.Cells.Item(138, 14).Value = -20780


# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 9262.373
$ws.Cells.Item(32, 9).Value = 5116.885
$ws.Cells.Item(32, 11).Value = 5116.885
$ws.Cells.Item(32, 13).Value = -4829.885

# Row 61
$ws.Cells.Item(61, 8).Value = 3724.2727
$ws.Cells.Item(61, 9).Value = 3619.111
$ws.Cells.Item(61, 11).Value = 3619.111
$ws.Cells.Item(61, 13).Value = -3407.111

# Row 132
$ws.Cells.Item(132, 8).Value = 4997.727
$ws.Cells.Item(132, 9).Value = 3663.3333
$ws.Cells.Item(132, 11).Value = 10989.9999
$ws.Cells.Item(132, 13).Value = -8459.999899999999

# Row 136
$ws.Cells.Item(136, 8).Value = 3724.2727
$ws.Cells.Item(136, 9).Value = 3619.111
$ws.Cells.Item(136, 11).Value = 10857.333
$ws.Cells.Item(136, 13).Value = -8307.332999999999


# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Cells.Item(22, 8).Value = 2663.7
$ws.Cells.Item(22, 9).Value = 2854.625
$ws.Cells.Item(22, 10).Value = 1900
$ws.Cells.Item(22, 11).Value = 2854.625
$ws.Cells.Item(22, 12).Value = 1900
$ws.Cells.Item(22, 13).Value = -2681.625
$ws.Cells.Item(22, 14).Value = -2246

# Row 86
$ws.Cells.Item(86, 8).Value = 5885006.5
$ws.Cells.Item(86, 10).Value = 1417.7778
$ws.Cells.Item(86, 12).Value = 1417.7778
$ws.Cells.Item(86, 14).Value = -3663.7778

# Row 89
$ws.Cells.Item(89, 8).Value = 5885006.5
$ws.Cells.Item(89, 10).Value = 1417.7778
$ws.Cells.Item(89, 12).Value = 7088.889
$ws.Cells.Item(89, 14).Value = -18320.889

# Row 99
$ws.Cells.Item(99, 8).Value = 5294457
$ws.Cells.Item(99, 10).Value = 4384.5713
$ws.Cells.Item(99, 12).Value = 4384.5713
$ws.Cells.Item(99, 14).Value = -7380.5713

# Row 109
$ws.Cells.Item(109, 14).ClearContents()
$ws.Cells.Item(109, 8).Value = 0
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 12).Value = 0

# Row 134
$ws.Cells.Item(134, 8).Value = 3226.45
$ws.Cells.Item(134, 9).Value = 1589.5312
$ws.Cells.Item(134, 10).Value = 9774.125
$ws.Cells.Item(134, 11).Value = 4768.5936
$ws.Cells.Item(134, 12).Value = 29322.375
$ws.Cells.Item(134, 13).Value = -2233.5936
$ws.Cells.Item(134, 14).Value = -34392.375


# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 23
$ws.Cells.Item(23, 13).ClearContents()
$ws.Cells.Item(23, 8).Value = 50010
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 10).Value = 50010
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 12).Value = 50010
$ws.Cells.Item(23, 14).Value = -50490

# Row 27
$ws.Cells.Item(27, 13).ClearContents()
$ws.Cells.Item(27, 8).Value = 50010
$ws.Cells.Item(27, 9).Value = 0
$ws.Cells.Item(27, 10).Value = 50010
$ws.Cells.Item(27, 11).Value = 0
$ws.Cells.Item(27, 12).Value = 50010
$ws.Cells.Item(27, 14).Value = -50394

# Row 56
$ws.Cells.Item(56, 8).Value = 12761.111
$ws.Cells.Item(56, 9).Value = 9356.375
$ws.Cells.Item(56, 10).Value = 39999
$ws.Cells.Item(56, 11).Value = 9356.375
$ws.Cells.Item(56, 12).Value = 39999
$ws.Cells.Item(56, 13).Value = -8511.375
$ws.Cells.Item(56, 14).Value = -41689

# Row 86
$ws.Cells.Item(86, 8).Value = 13849.85
$ws.Cells.Item(86, 9).Value = 10300
$ws.Cells.Item(86, 11).Value = 10300
$ws.Cells.Item(86, 13).Value = -9177

# Row 89
$ws.Cells.Item(89, 8).Value = 13849.85
$ws.Cells.Item(89, 9).Value = 10300
$ws.Cells.Item(89, 11).Value = 51500
$ws.Cells.Item(89, 13).Value = -45884

# Row 99
$ws.Cells.Item(99, 8).Value = 3822
$ws.Cells.Item(99, 9).Value = 3454.125
$ws.Cells.Item(99, 11).Value = 3454.125
$ws.Cells.Item(99, 13).Value = -1956.125

# Row 107
$ws.Cells.Item(107, 8).Value = 1919.2333
$ws.Cells.Item(107, 9).Value = 1543.4166
$ws.Cells.Item(107, 10).Value = 3422.5
$ws.Cells.Item(107, 11).Value = 1543.4166
$ws.Cells.Item(107, 12).Value = 3422.5
$ws.Cells.Item(107, 13).Value = 376.5834
$ws.Cells.Item(107, 14).Value = -7262.5

# Row 109
$ws.Cells.Item(109, 8).Value = 19307.857
$ws.Cells.Item(109, 10).Value = 19307.857
$ws.Cells.Item(109, 12).Value = 19307.857
$ws.Cells.Item(109, 14).Value = -21387.857

# Row 126
$ws.Cells.Item(126, 8).Value = 3822
$ws.Cells.Item(126, 9).Value = 3454.125
$ws.Cells.Item(126, 11).Value = 10362.375
$ws.Cells.Item(126, 13).Value = -7892.375

# Row 134
$ws.Cells.Item(134, 8).Value = 2272.016
$ws.Cells.Item(134, 9).Value = 1536.7646
$ws.Cells.Item(134, 11).Value = 4610.293799999999
$ws.Cells.Item(134, 13).Value = -2075.293799999999


# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Cells.Item(2, 8).Value = 218.62791
$ws.Cells.Item(2, 9).Value = 108.65
$ws.Cells.Item(2, 11).Value = 651.9000000000001
$ws.Cells.Item(2, 13).Value = -538.9000000000001

# Row 4
$ws.Cells.Item(4, 8).Value = 10812309
$ws.Cells.Item(4, 9).Value = 11892552
$ws.Cells.Item(4, 10).Value = 550000
$ws.Cells.Item(4, 11).Value = 35677656
$ws.Cells.Item(4, 12).Value = 1650000
$ws.Cells.Item(4, 13).Value = -35677544
$ws.Cells.Item(4, 14).Value = -1650224


# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Cells.Item(5, 8).Value = 3500
$ws.Cells.Item(5, 10).Value = 3500
$ws.Cells.Item(5, 12).Value = 3500
$ws.Cells.Item(5, 14).Value = -3724

# Row 20
$ws.Cells.Item(20, 8).Value = 10500
$ws.Cells.Item(20, 9).Value = 10500
$ws.Cells.Item(20, 11).Value = 10500
$ws.Cells.Item(20, 13).Value = -10255

# Row 132
$ws.Cells.Item(132, 8).Value = 4105.1763
$ws.Cells.Item(132, 9).Value = 3680.3
$ws.Cells.Item(132, 10).Value = 4712.143
$ws.Cells.Item(132, 11).Value = 11040.9
$ws.Cells.Item(132, 12).Value = 14136.429
$ws.Cells.Item(132, 13).Value = -8510.900000000001
$ws.Cells.Item(132, 14).Value = -19196.429


# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 11524.75
$ws.Cells.Item(7, 10).Value = 13333
$ws.Cells.Item(7, 12).Value = 13333
$ws.Cells.Item(7, 14).Value = -13557

# Row 9
$ws.Cells.Item(9, 8).Value = 1795.1428
$ws.Cells.Item(9, 9).Value = 1203.3334
$ws.Cells.Item(9, 11).Value = 1203.3334
$ws.Cells.Item(9, 13).Value = -979.3334

# Row 16
$ws.Cells.Item(16, 8).Value = 1780
$ws.Cells.Item(16, 10).Value = 5000
$ws.Cells.Item(16, 12).Value = 5000
$ws.Cells.Item(16, 14).Value = -5340

# Row 22
$ws.Cells.Item(22, 8).Value = 75696.836
$ws.Cells.Item(22, 10).Value = 2500
$ws.Cells.Item(22, 12).Value = 2500
$ws.Cells.Item(22, 14).Value = -3090

# Row 27
$ws.Cells.Item(27, 8).Value = 75696.836
$ws.Cells.Item(27, 10).Value = 2500
$ws.Cells.Item(27, 12).Value = 2500
$ws.Cells.Item(27, 14).Value = -2714

# Row 126
$ws.Cells.Item(126, 8).Value = 11524.75
$ws.Cells.Item(126, 10).Value = 13333
$ws.Cells.Item(126, 12).Value = 39999
$ws.Cells.Item(126, 14).Value = -44939


# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Cells.Item(81, 8).Value = 9263431
$ws.Cells.Item(81, 9).Value = 18520374
$ws.Cells.Item(81, 10).Value = 6487.4443
$ws.Cells.Item(81, 11).Value = 37040748
$ws.Cells.Item(81, 12).Value = 12974.8886
$ws.Cells.Item(81, 13).Value = -37039687
$ws.Cells.Item(81, 14).Value = -15096.8886

# Row 84
$ws.Cells.Item(84, 8).Value = 9263431
$ws.Cells.Item(84, 9).Value = 18520374
$ws.Cells.Item(84, 10).Value = 6487.4443
$ws.Cells.Item(84, 11).Value = 185203740
$ws.Cells.Item(84, 12).Value = 64874.443
$ws.Cells.Item(84, 13).Value = -185198436
$ws.Cells.Item(84, 14).Value = -75482.443

# Row 126
$ws.Cells.Item(126, 8).Value = 2423.818
$ws.Cells.Item(126, 9).Value = 2583
$ws.Cells.Item(126, 10).Value = 1999.3334
$ws.Cells.Item(126, 11).Value = 7749
$ws.Cells.Item(126, 12).Value = 5998.0002
$ws.Cells.Item(126, 13).Value = -5279
$ws.Cells.Item(126, 14).Value = -10938.0002

